$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "If an email is saved to drafts first then sent, the
#    attachment won't be visible for download" gets " (Fixed)" appended,
#    with "(Fixed)" in bold red.
# ---------------------------------------------------------------------------

# Locate the paragraph by its current (unique) text.
$targetText = "If an email is saved to drafts first then sent, the attachment " `
    + [char]0x2019 + "t be visible for download"

$issuesPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq $targetText) {
        $issuesPara = $i
        break
    }
}

# Append " (Fixed)" as a single chunk right after the existing text so the
# inserted text inherits the paragraph's bold/size formatting cleanly.
$p = $d.Paragraphs.Item($issuesPara)
$p.Range.InsertAfter(" (Fixed)")

# Colour just "(Fixed)" -- leaving the leading space uncoloured -- by finding
# each piece and setting its font colour individually. This creates three
# independent runs: "(", "Fixed" and ")".
$p = $d.Paragraphs.Item($issuesPara)
$rOpen = $p.Range
$rOpen.Find.Execute("(", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rOpen.Font.Color = 255

$p = $d.Paragraphs.Item($issuesPara)
$rWord = $p.Range
$rWord.Find.Execute("Fixed", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rWord.Font.Color = 255

$p = $d.Paragraphs.Item($issuesPara)
$rClose = $p.Range
$rClose.Find.Execute(")", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rClose.Font.Color = 255

# Finally, split the leading space off from the original "...download" run
# so it becomes its own run (matching how the document was actually edited)
# instead of being silently merged back into the preceding run.
$p = $d.Paragraphs.Item($issuesPara)
$rSpace = $p.Range
$rSpace.Find.Execute("download ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rSpace.MoveStart(1, 8) | Out-Null
$rSpace.Font.Bold = $false
$rSpace.Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) Add a new bullet "There seems to be a max size upload limit" right
#    after "Show time format better in folders ".
# ---------------------------------------------------------------------------

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $lastPara.Range.InsertParagraphAfter()

$addedPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$addedPara.Range.InsertAfter("There seems to be a max size upload limit")
